# Planilha com ampliação de incertezas
# Expand the uncertainty ranges in the "params" sheet by switching several
# parameters from "Fixo" (fixed) to "Incerto" (uncertain), which widens
# their Min/Max (columns C/D) via the existing formulas that depend on
# column I. A couple of rows also get explicit tweaks to F (variance
# factor) and G (lower plausible bound).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# Row 6: widen the variance factor and make the lower bound a formula (1/15)
$ws.Range("F6").Value = 2
$ws.Range("G6").Formula = "=1/15"

# Row 27 & 29: explicit lower plausible bound (column G) updates
# (set before the I-column switch below so the MAX()/MIN() formulas in
# C/D recompute off the new G value once "Incerto" takes effect)
$ws.Range("G27").Value = 0.7
$ws.Range("G29").Value = 0.6

# Rows whose "Tipo" (column I) switches from "Fixo" to "Incerto"
$incertoRows = @(11, 15, 17, 20, 21, 24, 26, 27, 29, 30, 34, 35, 36, 40, 42, 43, 45, 46, 49)
foreach ($r in $incertoRows) {
    $ws.Range("I$r").Value = "Incerto"
}

# Reset the active selection back to A1 (as in the saved workbook)
$ws.Range("A1").Select()
